$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correction post-soutenance : 5 clusters kmeans
# Update column M (Kmeans) values for the affected rows
$kmeansUpdates = @{
    3 = 5
    4 = 4
    6 = 4
    7 = 3
    8 = 4
    9 = 3
    10 = 4
    11 = 5
    12 = 2
    13 = 4
    14 = 5
    15 = 4
    17 = 5
    18 = 4
    19 = 5
    21 = 5
    22 = 5
    23 = 5
    24 = 5
    26 = 2
    27 = 3
    28 = 3
    31 = 4
    32 = 2
    33 = 4
    34 = 4
    35 = 4
    36 = 5
    37 = 4
    38 = 5
    39 = 3
    40 = 5
    41 = 4
    42 = 5
    44 = 4
    45 = 3
    46 = 4
    47 = 5
    48 = 4
    49 = 4
    50 = 5
    51 = 5
    52 = 4
    53 = 2
    54 = 5
    55 = 3
    56 = 3
    57 = 3
    58 = 5
    59 = 2
    60 = 5
    62 = 3
    63 = 5
    64 = 5
    65 = 3
    66 = 5
    67 = 2
    70 = 5
    71 = 3
    72 = 4
    73 = 4
    74 = 5
    75 = 2
    76 = 5
    77 = 3
    78 = 4
    80 = 3
    81 = 5
    82 = 3
    83 = 5
    84 = 4
    85 = 3
    86 = 3
    87 = 4
    88 = 4
    89 = 5
    91 = 2
    93 = 3
    95 = 4
    96 = 5
    97 = 5
    98 = 3
    99 = 5
    100 = 4
    101 = 4
    103 = 5
    104 = 3
    105 = 5
    108 = 4
    109 = 4
    110 = 4
    112 = 3
    115 = 5
    116 = 5
    117 = 4
    119 = 2
    120 = 4
    121 = 2
    123 = 2
    124 = 4
    127 = 2
    128 = 5
    129 = 5
    132 = 4
    133 = 5
    134 = 5
    135 = 5
    136 = 3
    137 = 5
    139 = 4
    140 = 4
    143 = 4
    144 = 5
    145 = 4
    147 = 3
    149 = 4
    150 = 5
    151 = 3
    152 = 3
    153 = 5
    154 = 5
    157 = 5
    158 = 4
    159 = 3
    160 = 5
    161 = 5
    165 = 1
    166 = 3
    167 = 5
    168 = 2
    170 = 3
}

foreach ($row in $kmeansUpdates.Keys) {
    $ws.Cells.Item($row, 13).Value = $kmeansUpdates[$row]
}

Write-Output "Updated $($kmeansUpdates.Count) Kmeans (column M) values"
